# Updated cryptos list on Thu Apr 18 22:08:55 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price ("D" column) values must stay plain text exactly as scraped, even
# when they look numeric (e.g. "550.79"), so force text formatting while
# writing them, then restore the default "Normal" style so no stray
# number-format style is left behind on the cell.
$prices = @{
    "D2"  = "63.517.80"
    "D3"  = "3.065.56"
    "D5"  = "550.79"
    "D6"  = "140.97"
    "D8"  = "3.060.22"
    "D10" = "6.49"
    "D12" = "0.456"
    "D13" = "0.0000228"
    "D14" = "34.98"
    "D15" = "3.566.15"
    "D16" = "63.542.15"
    "D17" = "3.067.40"
    "D19" = "6.78"
    "D20" = "485.11"
    "D22" = "0.677"
    "D23" = "7.29"
    "D24" = "81.17"
    "D25" = "12.75"
    "D28" = "7.90"
    "D30" = "1.00"
    "D31" = "26.24"
    "D33" = "2.46"
    "D34" = "5.68"
    "D36" = "6.00"
    "D37" = "465.13"
    "D38" = "0.0824"
    "D39" = "0.0399"
    "D40" = "3.054.85"
    "D42" = "8.24"
    "D44" = "27.97"
    "D45" = "0.256"
    "D47" = "2.05"
    "D48" = "0.110"
    "D49" = "117.11"
    "D50" = "0.0₃0511"
    "D51" = "2.08"
}

foreach ($addr in $prices.Keys) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $prices[$addr]
    $ws.Range($addr).Style = "Normal"
}

# Volume(1h) ("E" column) values are padded with spaces around the sign
# and percent, so Excel already treats them as plain text.
$volumes = @{
    "E2"  = "  +3.73%  "
    "E3"  = "  +2.50%  "
    "E4"  = "  -0.01%  "
    "E5"  = "  +2.93%  "
    "E6"  = "  +4.61%  "
    "E7"  = "  -0.09%  "
    "E8"  = "  +2.51%  "
    "E9"  = "  +1.50%  "
    "E10" = "  +5.58%  "
    "E11" = "  +3.22%  "
    "E12" = "  +2.22%  "
    "E13" = "  +3.09%  "
    "E14" = "  +3.13%  "
    "E15" = "  +2.55%  "
    "E16" = "  +3.67%  "
    "E17" = "  +2.42%  "
    "E18" = "  -0.95%  "
    "E19" = "  +2.50%  "
    "E20" = "  +4.47%  "
    "E21" = "  +4.57%  "
    "E22" = "  +0.24%  "
    "E23" = "  +5.07%  "
    "E24" = "  +1.46%  "
    "E25" = "  +6.50%  "
    "E26" = "  +0.07%  "
    "E27" = "  +3.46%  "
    "E28" = "  +1.36%  "
    "E29" = "  +6.93%  "
    "E30" = "  -0.07%  "
    "E31" = "  +2.72%  "
    "E32" = "  +1.19%  "
    "E33" = "  +8.30%  "
    "E34" = "  +3.98%  "
    "E35" = "  +0.31%  "
    "E36" = "  +1.94%  "
    "E37" = "  +3.19%  "
    "E38" = "  +4.61%  "
    "E39" = "  +3.76%  "
    "E40" = "  -3.22%  "
    "E41" = "  +1.17%  "
    "E42" = "  +1.44%  "
    "E43" = "  +4.70%  "
    "E44" = "  +2.91%  "
    "E45" = "  +4.87%  "
    "E47" = "  +2.71%  "
    "E48" = "  +2.60%  "
    "E49" = "  -1.57%  "
    "E50" = "  +3.25%  "
    "E51" = "  +4.46%  "
}

foreach ($addr in $volumes.Keys) {
    $ws.Range($addr).Value = $volumes[$addr]
}

# Rows 49 and 50 swap coins: PEPE (previously row 49) moves to row 50,
# and Monero (previously row 50) moves to row 49.
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"

$ws.Range("B50").Value = "PEPE"
$ws.Range("C50").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
